$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "WOfVH549"
$ws.Range("B2").Value = 231011198
$ws.Range("C2").Value = "wxilozm91"
$ws.Range("D2").Value = "C3&d2%Gr"
$ws.Range("F2").Value = "fyNmpqfB"
$ws.Range("G2").Value = "nRKq"
